$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42:B55").ClearContents()

$ws.Cells.Item(2, 1).Value = 'Home_Last 5_Margin'
$ws.Cells.Item(2, 2).Value = 0.1112767007389389
$ws.Cells.Item(3, 1).Value = 'Prev_Rush'
$ws.Cells.Item(3, 2).Value = 0.09427479998943883
$ws.Cells.Item(4, 1).Value = 'Away _Last 5_ Margin'
$ws.Cells.Item(4, 2).Value = 0.07671970318667719
$ws.Cells.Item(5, 1).Value = 'A_Prev_Pass'
$ws.Cells.Item(5, 2).Value = 0.0680886434333324
$ws.Cells.Item(6, 1).Value = 'Week_TO_Margin_Interaction'
$ws.Cells.Item(6, 2).Value = 0.04508925949423884
$ws.Cells.Item(7, 1).Value = 'Prev_Spec'
$ws.Cells.Item(7, 2).Value = 0.04308770721675662
$ws.Cells.Item(8, 1).Value = 'Away_PrevYrW'
$ws.Cells.Item(8, 2).Value = 0.04291532262166163
$ws.Cells.Item(9, 1).Value = 'A_Prev_Rush'
$ws.Cells.Item(9, 2).Value = 0.04079322261163119
$ws.Cells.Item(10, 1).Value = 'Week_WinInteraction'
$ws.Cells.Item(10, 2).Value = 0.04037415091441079
$ws.Cells.Item(11, 1).Value = 'Prev_Tot_Def'
$ws.Cells.Item(11, 2).Value = 0.04026684438206787
$ws.Cells.Item(12, 1).Value = 'A_Prev_Tot_Def'
$ws.Cells.Item(12, 2).Value = 0.04017440186146815
$ws.Cells.Item(13, 1).Value = 'Prev_Pass'
$ws.Cells.Item(13, 2).Value = 0.0388112417699172
$ws.Cells.Item(14, 1).Value = 'A_Prev_Spec'
$ws.Cells.Item(14, 2).Value = 0.0384022732972765
$ws.Cells.Item(15, 1).Value = 'Away_TO_Margin_Interaction'
$ws.Cells.Item(15, 2).Value = 0.03635120387902429
$ws.Cells.Item(16, 1).Value = 'Away_Penalty_Yards'
$ws.Cells.Item(16, 2).Value = 0.03528162460422748
$ws.Cells.Item(17, 1).Value = 'Week_Third_Down'
$ws.Cells.Item(17, 2).Value = 0.03508914130488096
$ws.Cells.Item(18, 1).Value = 'Away_WinInteraction'
$ws.Cells.Item(18, 2).Value = 0.03481521523097079
$ws.Cells.Item(19, 1).Value = 'Away_Third_Down'
$ws.Cells.Item(19, 2).Value = 0.03402573776446908
$ws.Cells.Item(20, 1).Value = 'Home_PrevYrW'
$ws.Cells.Item(20, 2).Value = 0.01857825092209707
$ws.Cells.Item(21, 1).Value = 'H_Wins5'
$ws.Cells.Item(21, 2).Value = 0.01207379355938058
$ws.Cells.Item(22, 1).Value = 'A_Wins5'
$ws.Cells.Item(22, 2).Value = 0.01140159348374387
$ws.Cells.Item(23, 1).Value = 'Week'
$ws.Cells.Item(23, 2).Value = 0.00974243199044653
$ws.Cells.Item(24, 1).Value = 'AwayDiv_AFC North'
$ws.Cells.Item(24, 2).Value = 0.007059540185243753
$ws.Cells.Item(25, 1).Value = 'HomeDiv_NFC East'
$ws.Cells.Item(25, 2).Value = 0.006494758910149957
$ws.Cells.Item(26, 1).Value = 'Home_QBInjury'
$ws.Cells.Item(26, 2).Value = 0.004269137864614808
$ws.Cells.Item(27, 1).Value = 'HomeDiv_NFC West'
$ws.Cells.Item(27, 2).Value = 0.003838526969471777
$ws.Cells.Item(28, 1).Value = 'HomeDiv_AFC North'
$ws.Cells.Item(28, 2).Value = 0.003714904396126878
$ws.Cells.Item(29, 1).Value = 'AwayDiv_NFC South'
$ws.Cells.Item(29, 2).Value = 0.00361301408343365
$ws.Cells.Item(30, 1).Value = 'AwayDiv_AFC South'
$ws.Cells.Item(30, 2).Value = 0.003530893725743437
$ws.Cells.Item(31, 1).Value = 'AwayDiv_AFC West'
$ws.Cells.Item(31, 2).Value = 0.003225943223864184
$ws.Cells.Item(32, 1).Value = 'HomeDiv_NFC North'
$ws.Cells.Item(32, 2).Value = 0.00291349669610598
$ws.Cells.Item(33, 1).Value = 'HomeDiv_NFC South'
$ws.Cells.Item(33, 2).Value = 0.002590572972765159
$ws.Cells.Item(34, 1).Value = 'Away_QBInjury'
$ws.Cells.Item(34, 2).Value = 0.002509746243145366
$ws.Cells.Item(35, 1).Value = 'AwayDiv_NFC East'
$ws.Cells.Item(35, 2).Value = 0.002127674277777717
$ws.Cells.Item(36, 1).Value = 'AwayDiv_NFC West'
$ws.Cells.Item(36, 2).Value = 0.001915255323526914
$ws.Cells.Item(37, 1).Value = 'HomeDiv_AFC West'
$ws.Cells.Item(37, 2).Value = 0.001730694935832774
$ws.Cells.Item(38, 1).Value = 'AwayDiv_NFC North'
$ws.Cells.Item(38, 2).Value = 0.001484409597807534
$ws.Cells.Item(39, 1).Value = 'HomeDiv_AFC South'
$ws.Cells.Item(39, 2).Value = 0.001348166337333254
$ws.Cells.Item(40, 1).Value = 'HomeTeam'
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 1).Value = 'AwayTeam'
$ws.Cells.Item(41, 2).Value = 0
